$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.230.45"
$ws.Range("E2").Value = "'  -1.33%  "
$ws.Range("D3").Value = "'1.877.77"
$ws.Range("E3").Value = "'  -0.11%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("D5").Value = "'235.18"
$ws.Range("E5").Value = "'  -0.91%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("D7").Value = "'0.4847"
$ws.Range("E7").Value = "'  -0.50%  "
$ws.Range("D8").Value = "'0.2875"
$ws.Range("E8").Value = "'  -2.31%  "
$ws.Range("D9").Value = "'0.06566"
$ws.Range("E9").Value = "'  -2.01%  "
$ws.Range("D10").Value = "'1.885.54"
$ws.Range("E10").Value = "'  +0.33%  "
$ws.Range("D11").Value = "'16.69"
$ws.Range("E11").Value = "'  -0.15%  "
$ws.Range("D12").Value = "'0.07259"
$ws.Range("E12").Value = "'  -0.77%  "
$ws.Range("D13").Value = "'5.110"
$ws.Range("E13").Value = "'  +1.21%  "
$ws.Range("D14").Value = "'87.02"
$ws.Range("E14").Value = "'  -2.79%  "
$ws.Range("D15").Value = "'0.6532"
$ws.Range("E15").Value = "'  -1.00%  "
$ws.Range("D16").Value = "'30.221.37"
$ws.Range("E16").Value = "'  -1.19%  "
$ws.Range("B17").Value = "'Dai"
$ws.Range("C17").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "'  +0.04%  "
$ws.Range("D18").Value = "'13.25"
$ws.Range("E18").Value = "'  +0.99%  "
$ws.Range("B19").Value = "'ShibaInu"
$ws.Range("C19").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007715"
$ws.Range("E19").Value = "'  -1.65%  "
$ws.Range("B20").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "'2.131.83"
$ws.Range("E20").Value = "'  +0.46%  "
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.294"
$ws.Range("E21").Value = "'  +11.28%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.09%  "
$ws.Range("D23").Value = "'192.01"
$ws.Range("E23").Value = "'  -10.38%  "
$ws.Range("D24").Value = "'6.097"
$ws.Range("E24").Value = "'  -0.87%  "
$ws.Range("D25").Value = "'9.312"
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("D26").Value = "'160.85"
$ws.Range("E26").Value = "'  +2.29%  "
$ws.Range("D27").Value = "'17.94"
$ws.Range("E27").Value = "'  -5.50%  "
$ws.Range("D28").Value = "'1.902"
$ws.Range("E28").Value = "'  +3.11%  "
$ws.Range("D29").Value = "'1.441"
$ws.Range("E29").Value = "'  +2.14%  "
$ws.Range("D30").Value = "'4.243"
$ws.Range("E30").Value = "'  -0.68%  "
$ws.Range("D31").Value = "'0.09115"
$ws.Range("E31").Value = "'  -0.07%  "
$ws.Range("D32").Value = "'4.037"
$ws.Range("E32").Value = "'  +1.57%  "
$ws.Range("D33").Value = "'0.05115"
$ws.Range("E33").Value = "'  -0.23%  "
$ws.Range("D34").Value = "'0.7237"
$ws.Range("E34").Value = "'  -2.32%  "
$ws.Range("D35").Value = "'1.093"
$ws.Range("E35").Value = "'  +0.47%  "
$ws.Range("D36").Value = "'2.700"
$ws.Range("E36").Value = "'  +1.23%  "
$ws.Range("D37").Value = "'0.01794"
$ws.Range("E37").Value = "'  -1.51%  "
$ws.Range("D38").Value = "'2.639"
$ws.Range("E38").Value = "'  -1.25%  "
$ws.Range("D39").Value = "'0.9144"
$ws.Range("E39").Value = "'  -0.41%  "
$ws.Range("D40").Value = "'2.038"
$ws.Range("E40").Value = "'  -1.67%  "
$ws.Range("B41").Value = "'Quant"
$ws.Range("C41").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'105.75"
$ws.Range("E41").Value = "'  -0.20%  "
$ws.Range("B42").Value = "'TheSandbox"
$ws.Range("C42").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.4286"
$ws.Range("E42").Value = "'  -3.62%  "
$ws.Range("D43").Value = "'5.805"
$ws.Range("E43").Value = "'  +0.90%  "
$ws.Range("D44").Value = "'0.9990"
$ws.Range("E44").Value = "'  +0.29%  "
$ws.Range("B45").Value = "'Aave"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'66.38"
$ws.Range("E45").Value = "'  +15.68%  "
$ws.Range("B46").Value = "'Aptos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.414"
$ws.Range("E46").Value = "'  -0.99%  "
$ws.Range("B47").Value = "'Algorand"
$ws.Range("C47").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1320"
$ws.Range("E47").Value = "'  -1.61%  "
$ws.Range("D48").Value = "'8.976"
$ws.Range("E48").Value = "'  +2.90%  "
$ws.Range("B49").Value = "'Cronos"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05773"
$ws.Range("E49").Value = "'  -1.43%  "
$ws.Range("B50").Value = "'Elrond"
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'33.86"
$ws.Range("E50").Value = "'  +1.18%  "
$ws.Range("D51").Value = "'0.3833"
$ws.Range("E51").Value = "'  -4.54%  "
